# The original file had G6 = "yes" (response to the review question in column G).
# The uploaded edit simply corrects/changes that answer to "no".
# All the other index shuffling visible in the raw OOXML diff is just a side
# effect of Excel/LibreOffice re-serializing the shared-strings table after
# this single cell edit - no other cell's displayed value actually changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "no"

# Reflect the resulting selection/active-cell position recorded in the diff
# (the cursor ends up one row below the edited cell, on G7).
$ws.Range("G7").Select()
